$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("P_req")
$ws2 = $wb.Worksheets.Item("Q_req")

# Build the new column B values for rows 2-52 (51 values):
#  - rows 2-12  (11 rows) -> -1500000
#  - rows 13-52 (40 rows) -> 200000
$values = @()
for ($i = 0; $i -lt 11; $i++) { $values += -1500000 }
for ($i = 0; $i -lt 40; $i++) { $values += 200000 }

for ($i = 0; $i -lt 51; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $values[$i]
    $ws2.Cells.Item($row, 2).Value = $values[$i]
}

# Update the view/selection on each sheet to match the new selection
# (B2:B32 active at B2) and reset scroll position to the top (A1).
$ws1.Activate()
$ws1.Range("B2:B32").Select()
$excel.ActiveWindow.ScrollRow = 1

$ws2.Activate()
$ws2.Range("B2:B32").Select()
$excel.ActiveWindow.ScrollRow = 1

$ws1.Activate()
